# Updated Indonesia files compatible with v3.3.1
# Applies the "PPEIdtIL" workbook edits:
#  - About sheet: insert new explanatory paragraph, push old note text down
#  - PPEIdtIL sheet: new wrapped header for column A, retitled/swapped
#    column headers, and updated percentages (0.1 -> 0.02)

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("PPEIdtIL")

# ---------------------------------------------------------------------
# About sheet
# ---------------------------------------------------------------------

# New paragraph explaining the variable (inserted above the old note,
# pushing it from rows 10-14 down to rows 18-22).
$wsAbout.Range("A10").Value = 'This variable reflects improvement in efficiency components selected by consumers due'
$wsAbout.Range("A11").Value = 'to improved labeling. The labeling influences consumers who are buying appliances of all'
$wsAbout.Range("A12").Value = 'Quality levels, so it''s represented as a simple percentage increase in the efficiency of'
$wsAbout.Range("A13").Value = 'components sold (at all quality levels). If Quality Levels are defined based on'
$wsAbout.Range("A14").Value = 'particular efficiency thresholds, this may mean the number of square feet served by'
$wsAbout.Range("A15").Value = 'components of a given quality level will not be accurate. It''s just a question of the meaning'
$wsAbout.Range("A16").Value = 'of the labels given to each Quality Level.'

# Row 17 intentionally left blank (spacer row).

# The original ACEEE sourcing note, now relocated to rows 18-22.
$wsAbout.Range("A18").Value = 'The ACEEE study focused on "appliance" labeling.  We use the same percentage for labeling'
$wsAbout.Range("A19").Value = 'of heating equipment, as well as cooling and ventilation equipment, because they are'
$wsAbout.Range("A20").Value = 'similar (e.g. machines one buys in a store, which could readily bear labels, with similar'
$wsAbout.Range("A21").Value = 'costs and lifetimes as other major appliances).  We similarly assume the same rate'
$wsAbout.Range("A22").Value = 'applies to commercial and residential buildings.'

# ---------------------------------------------------------------------
# PPEIdtIL sheet
# ---------------------------------------------------------------------

# Row 1 header relabel: A1 becomes a wrapped descriptive title, and the
# region columns are reordered (Urban Residential, Rural Residential,
# Commercial).
$wsData.Range("A1").Value = "Efficiency Improvement by Building Component (dimensionless)"
$wsData.Range("A1").WrapText = $true
$wsData.Rows.Item(1).RowHeight = 45

$wsData.Range("B1").Value = "Urban Residential"
$wsData.Range("C1").Value = "Rural Residential"
$wsData.Range("D1").Value = "Commercial"

# Updated percentages: heating, cooling & ventilation, and appliances
# move from 10% to 2%. Envelope, lighting, and other component remain 0.
$wsData.Range("B2").Value = 0.02
$wsData.Range("D2").Value = 0.02

$wsData.Range("B3").Value = 0.02
$wsData.Range("D3").Value = 0.02

$wsData.Range("B6").Value = 0.02
$wsData.Range("D6").Value = 0.02

# ---------------------------------------------------------------------
# View state: PPEIdtIL becomes the active/selected tab (cell L3),
# About keeps a lingering selection at E35 from when it was last active.
# ---------------------------------------------------------------------

$wsAbout.Activate()
[void]$wsAbout.Range("E35").Select()

$wsData.Activate()
[void]$wsData.Range("L3").Select()
